# Add a new "Save" column (H) to the s_vals sheet, matching the
# existing header formatting (bold font, thin border, centered/top
# aligned) and fill the data rows with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, styled like the other header cells (copy format from G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data values.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
